$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.70015025138855
$ws.Range("B1").Value = 3.102221250534058
$ws.Range("C1").Value = 2.614607810974121
$ws.Range("D1").Value = 2.056376695632935
$ws.Range("E1").Value = 1.28082013130188
